# Applies the "Alvearie FHIR IG" StructureDefinition metadata refresh:
#   - Version bump 5.0.0 -> 6.0.0
#   - Date refresh
#   - Publisher's Contact/"No display for ContactDetail" rows replaced with
#     a real Publisher value ("Alvearie Team") plus a new Jurisdiction row
#   - The differential's root Extension row's Short/Definition now mirror
#     the StructureDefinition's own Title/Description

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version
$meta.Range("B3").Value = "6.0.0"

# Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$meta.Range("B9").Value = "Alvearie Team"

# The old sheet had two duplicate "Contact" / "No display for ContactDetail"
# rows right after Publisher - remove both.
$meta.Rows.Item(10).Delete()
$meta.Rows.Item(10).Delete()

# ... and replace them with a single new "Jurisdiction" row, formatted like
# the rest of the table.
$meta.Rows.Item(10).Insert()
$meta.Range("A9:B9").Copy()
$meta.Range("A10:B10").PasteSpecial(-4122)
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# The "Elements" differential sheet: the top-level Extension row's Short /
# Definition columns now reuse the StructureDefinition's Title/Description
# instead of the generic "Extension" / "An Extension" placeholders.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Utilization Review"
$elements.Range("L2").Value = "Customer-specific code for the type of utilization review"
